# Update countries & provincias Spain
# Applies the data refresh captured in the commit: new totals for several
# countries, the swap of "Montserrat" / "Islas Malvinas" row data, and the
# refreshed "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 07:41"

# --- India (row 6) -------------------------------------------------------
$ws.Range("B6").Value = 2461542
$ws.Range("C6").Value = 1929
$ws.Range("D6").Value = 1751846
$ws.Range("E6").Value = 661543
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 48153

# --- Kazajistan (row 29) --------------------------------------------------
$ws.Range("D29").Value = 78633
$ws.Range("E29").Value = 21946

# --- Israel (row 33) -------------------------------------------------------
$ws.Range("B33").Value = 90283
$ws.Range("C33").Value = 461
$ws.Range("D33").Value = 66143
$ws.Range("E33").Value = 23489

# --- Kirguistan (row 55) ----------------------------------------------------
$ws.Range("B55").Value = 41373
$ws.Range("C55").Value = 304
$ws.Range("D55").Value = 33592
$ws.Range("E55").Value = 6290
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 1491

# --- Uzbekistan (row 62) ----------------------------------------------------
$ws.Range("B62").Value = 33561
$ws.Range("C62").Value = 238
$ws.Range("E62").Value = 6130
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 218

# --- Bulgaria (row 81) -------------------------------------------------------
$ws.Range("B81").Value = 14069
$ws.Range("D81").Value = 8901
$ws.Range("E81").Value = 4684
$ws.Range("H81").Value = 484

# --- Haiti (row 95) -----------------------------------------------------------
$ws.Range("B95").Value = 7810
$ws.Range("C95").Value = 29
$ws.Range("E95").Value = 2495

# --- Tailandia (row 117) -------------------------------------------------------
$ws.Range("B117").Value = 3376
$ws.Range("C117").Value = 17
$ws.Range("D117").Value = 3173
$ws.Range("E117").Value = 145

# --- Swap Islas Malvinas / Montserrat (rows 213-214) --------------------------
# Row 213 previously held "Islas Malvinas" (D=13, H=0); it now holds
# "Montserrat" data (D=12, H=1). Row 214 previously held "Montserrat"
# (D=12, H=1); it now holds "Islas Malvinas" data (D=13, H=0).
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
